$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before the current row 7 ("Think about passing...")
# This pushes the existing rows 7,9,10,12,14,17 down to 12,14,15,17,19,22
$ws.Rows("7:11").Insert()

# Introduce the brand-new shared strings in the exact order they were
# originally authored so the resulting shared-string table indices line up:
#   17 Test suit for SliderAdjuster
#   18 Add different Sheets in Factory
#   19 Move Factory to LocalMock directory and create Abstract class for factory service
#   20 Fix constructor to use DI in SheetSearchCriteria
#   21 Use slider a a toggle button instead of button
#   22 Post sales
#   23 Editing of Sheets by central user
#   24 Think about Accordion

$ws.Range("B24").Value = "Test suit for SliderAdjuster"
$ws.Range("B7").Value = "Add different Sheets in Factory"
$ws.Range("B8").Value = "Move Factory to LocalMock directory and create Abstract class for factory service"
$ws.Range("B9").Value = "Fix constructor to use DI in SheetSearchCriteria"
$ws.Range("B10").Value = "Use slider a a toggle button instead of button"
$ws.Range("B26").Value = "Post sales"
$ws.Range("B28").Value = "Editing of Sheets by central user"
$ws.Range("B11").Value = "Think about Accordion"

# Fill in the Status column (reusing existing "Open"/"Closed" shared strings)
$ws.Range("C7").Value = "Open"
$ws.Range("C8").Value = "Closed"
$ws.Range("C9").Value = "Closed"
$ws.Range("C10").Value = "Open"
$ws.Range("C11").Value = "Open"
$ws.Range("C24").Value = "Open"
$ws.Range("C26").Value = "Open"

# Restore the originally-selected cell
$ws.Range("B14").Select()
